$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59, shifting existing rows 59:87 down to 60:88
# (row formatting/styles of the old row 59, e.g. the date style on column D,
# carry down with the shift).
$ws.Range("A59:R59").Insert()

# Populate the newly inserted row with this week's data point for Cilantro.
$ws.Cells.Item(59, 1).Value = 1
$ws.Cells.Item(59, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(59, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(59, 4).Value = 44825
$ws.Cells.Item(59, 5).Value = 15
$ws.Cells.Item(59, 6).Value = 100112040
$ws.Cells.Item(59, 7).Value = "Cilantro"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 270
$ws.Cells.Item(59, 11).Value = 900
$ws.Cells.Item(59, 12).Value = 1000
$ws.Cells.Item(59, 13).Value = 950
$ws.Cells.Item(59, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(59, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(59, 16).Value = 475
$ws.Cells.Item(59, 17).Value = 2
$ws.Cells.Item(59, 18).Value = "Hortaliza"
